{"js": "// Title run \"C SK9822 Module Documentation\" is split into two runs with the\n// same formatting:\n//   run 1: \"Appendix B.\"\n//   run 2: \" SK9822 Module Documentation\"\n// (net effect: \"C \" -> \"Appendix B. \" while keeping \" SK9822 Module\n// Documentation\" as its own run, matching the authored OOXML diff.)\nconst body = context.document.body;\nconst results = body.search(\"C SK9822 Module Documentation\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the title text to update.\");\n}\n\nconst hit = results.items[0];\n\n// Use insertOoxml (Flat OPC) to replace the matched range's contents with two\n// runs sharing the original title formatting (CMU Serif, 36 half-points),\n// mirroring the two <w:r> elements the diff introduces.\nconst flatOpc = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"CMU Serif\" w:hAnsi=\"CMU Serif\" w:cs=\"CMU Serif\"/>\n                <w:sz w:val=\"36\"/>\n                <w:szCs w:val=\"36\"/>\n              </w:rPr>\n              <w:t>Appendix B.</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"CMU Serif\" w:hAnsi=\"CMU Serif\" w:cs=\"CMU Serif\"/>\n                <w:sz w:val=\"36\"/>\n                <w:szCs w:val=\"36\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\"> SK9822 Module Documentation</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nhit.insertOoxml(flatOpc, \"Replace\");\nawait context.sync();\n", "ps1": "# Title run \"C SK9822 Module Documentation\" is split into two runs with the\n# same formatting:\n#   run 1: \"Appendix B.\"\n#   run 2: \" SK9822 Module Documentation\"\n# (net effect: \"C \" -> \"Appendix B. \" while keeping \" SK9822 Module\n# Documentation\" as its own run, matching the authored OOXML diff.)\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$found = $range.Find.Execute(\"C SK9822 Module Documentation\")\n\nif (-not $found) {\n    throw \"Could not find the title text to update.\"\n}\n\n# Replace the matched range's contents with two runs sharing the original\n# title formatting (CMU Serif, 36 half-points), mirroring the two <w:r>\n# elements the diff introduces.\n$flatOpc = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"CMU Serif\" w:hAnsi=\"CMU Serif\" w:cs=\"CMU Serif\"/>\n                <w:sz w:val=\"36\"/>\n                <w:szCs w:val=\"36\"/>\n              </w:rPr>\n              <w:t>Appendix B.</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"CMU Serif\" w:hAnsi=\"CMU Serif\" w:cs=\"CMU Serif\"/>\n                <w:sz w:val=\"36\"/>\n                <w:szCs w:val=\"36\"/>\n              </w:rPr>\n              <w:t xml:space=\"preserve\"> SK9822 Module Documentation</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$range.InsertXML($flatOpc)\n"}
